$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 17: date, hours, task description
$ws.Cells.Item(17, 1).Value = 43991
$ws.Cells.Item(17, 2).Value = 1.75
$ws.Cells.Item(17, 3).Value = "Figure of share of sales: task 3.4 options for weighted survival function and checking S(1)"

# Match formatting of the row above (date style for A, text style for C)
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(17, 1).PasteSpecial(-4122)

$ws.Cells.Item(15, 3).Copy()
$ws.Cells.Item(17, 3).PasteSpecial(-4122)

$ws.Rows.Item(17).RowHeight = $ws.Rows.Item(15).RowHeight

$excel.CutCopyMode = 0

# Update selection to the newly added task cell
$ws.Range("C17").Select()
